$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change the address header cell from Arabic "   العنوان" to English "address"
$ws.Range("B1").Value = "address"

# Update the selected cell to F3
$ws.Range("F3").Select()
